$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(4, 6, 9, 10, 11, 13)
foreach ($r in $rows) {
    $ws.Range("G$r").Value = "lipid/free"
}
